# Add files via upload
# - Rename the "TP" column header to "AP" on both the for_A and for_B sheets
#   (new shared strings) and replace that column's data with the new AP
#   values.
# - Switch the active sheet from for_A to for_B, and update each sheet's
#   remembered selection/cursor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# for_A sheet: column I was "TP", becomes "AP" with new values
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("for_A")

$wsA.Range("I1").Value = "AP"

$wsA.Range("I2").Value = 0.29194999999999999
$wsA.Range("I3").Value = 0.10818
$wsA.Range("I4").Value = 0.46046999999999999
$wsA.Range("I5").Value = 0.090100000000000013
$wsA.Range("I6").Value = 0.10872
$wsA.Range("I7").Value = 0.069379999999999997
$wsA.Range("I8").Value = 0.19462000000000002
$wsA.Range("I9").Value = 0.42591999999999997
$wsA.Range("I10").Value = 0.12146999999999999
$wsA.Range("I11").Value = 0.20673
$wsA.Range("I12").Value = 0.13450000000000001
$wsA.Range("I13").Value = 0.18686999999999998
$wsA.Range("I14").Value = 0.21635000000000001
$wsA.Range("I15").Value = 0.11924999999999999
$wsA.Range("I16").Value = 0.12483

# ---------------------------------------------------------------------
# for_B sheet: column F was "TP", becomes "AP" with new values
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("for_B")

$wsB.Range("F1").Value = "AP"

$wsB.Range("F2").Value = 0.29194999999999999
$wsB.Range("F3").Value = 0.10818
$wsB.Range("F4").Value = 0.46046999999999999
$wsB.Range("F5").Value = 0.090100000000000013
$wsB.Range("F6").Value = 0.10872
$wsB.Range("F7").Value = 0.069379999999999997
$wsB.Range("F8").Value = 0.20673
$wsB.Range("F9").Value = 0.13450000000000001
$wsB.Range("F10").Value = 0.18686999999999998
$wsB.Range("F11").Value = 0.21635000000000001
$wsB.Range("F12").Value = 0.11924999999999999
$wsB.Range("F13").Value = 0.12483

# ---------------------------------------------------------------------
# View state: make for_A's selection a multi-area one (best effort — the
# interop Select() collapses to the first area), move the cursor, and
# make for_B the active/selected sheet.
# ---------------------------------------------------------------------
$wsA.Activate()
$areaA1 = $wsA.Range("I2:I7")
$areaA2 = $wsA.Range("I11:I16")
$unionA = $excel.Union($areaA1, $areaA2)
$unionA.Select() | Out-Null

$wsB.Activate()
$wsB.Range("E20").Select() | Out-Null
